$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 610.6667
$arr[0,1] = 610.6667
$arr[0,2] = 0
$arr[0,3] = 5496.0003
$arr[0,4] = 0
$arr[0,5] = -3036.0003
$ws.Range("H125:M125").Value = $arr

$ws.Range("H126").Value = 39260
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 39260
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 39260
$ws.Range("N126").Value = -49140

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2001.75
$arr[0,1] = 398.5
$arr[0,2] = 2536.1667
$arr[0,3] = 1195.5
$arr[0,4] = 7608.500100000001
$arr[0,5] = 3764.5
$arr[0,6] = -17528.5001
$ws.Range("H127:N127").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$ws.Range("H128:L128").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1234.6
$arr[0,1] = 1103.5
$arr[0,2] = 1259.5714
$arr[0,3] = 3310.5
$arr[0,4] = 3778.7142
$arr[0,5] = 1689.5
$arr[0,6] = -13778.7142
$ws.Range("H129:N129").Value = $arr

$ws.Range("H130").Value = 54980
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 54980
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 54980
$ws.Range("N130").Value = -65020

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1807.48
$arr[0,1] = 820.3684
$arr[0,2] = 4933.3335
$arr[0,3] = 2461.1052
$arr[0,4] = 14800.0005
$arr[0,5] = 2578.8948
$arr[0,6] = -24880.0005
$ws.Range("H131:N131").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 45107256
$arr[0,1] = 45107256
$arr[0,2] = 0
$arr[0,3] = 135321768
$arr[0,4] = 0
$arr[0,5] = -135319238
$ws.Range("H132:M132").Value = $arr

$ws.Range("H133").Value = 24162
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 24162
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 24162
$ws.Range("N133").Value = -34282

$ws.Range("H134").Value = 47500
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 47500
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 47500
$ws.Range("N134").Value = -57640

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4044
$arr[0,1] = 5068.8
$arr[0,2] = 2073.2307
$arr[0,3] = 45619.2
$arr[0,4] = 18659.0763
$arr[0,5] = -43084.2
$arr[0,6] = -23729.0763
$ws.Range("H135:N135").Value = $arr

$ws.Range("H136").Value = 20000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 20000
$ws.Range("N136").Value = -30200

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 9611.15
$arr[0,1] = 860.8276
$arr[0,2] = 17796.936
$arr[0,3] = 2582.4828
$arr[0,4] = 53390.808
$arr[0,5] = -32.48279999999977
$arr[0,6] = -58490.808
$ws.Range("H137:N137").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2647.442
$arr[0,1] = 1678.8
$arr[0,2] = 2940.9697
$arr[0,3] = 5036.4
$arr[0,4] = 8822.909100000001
$arr[0,5] = 103.6000000000004
$arr[0,6] = -19102.9091
$ws.Range("H138:N138").Value = $arr

$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 30000
$ws.Range("N139").Value = -40280

$ws.Range("H140").Value = 28000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 28000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 28000
$ws.Range("N140").Value = -38360

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3908.5715
$arr[0,1] = 3027.6086
$arr[0,2] = 7961
$arr[0,3] = 9082.825800000001
$arr[0,4] = 23883
$arr[0,5] = -3902.825800000001
$arr[0,6] = -34243
$ws.Range("H141:N141").Value = $arr

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 328364.5
$ws.Range("I32").Value = 2280.7092
$ws.Range("K32").Value = 2280.7092
$ws.Range("M32").Value = -1993.7092

$ws.Range("H88").Value = 15228.546
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 15228.546
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 15228.546
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -16040.546

$ws.Range("H91").Value = 15228.546
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 15228.546
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 15228.546
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -18036.546

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2360081
$arr[0,1] = 3677200.5
$arr[0,2] = 3130.3157
$arr[0,3] = 11031601.5
$arr[0,4] = 9390.947100000001
$arr[0,5] = -11029071.5
$arr[0,6] = -14450.9471
$ws.Range("H132:N132").Value = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1791.4286
$arr[0,1] = 1801.6
$arr[0,2] = 1766
$arr[0,3] = 1801.6
$arr[0,4] = 1766
$arr[0,5] = -678.5999999999999
$arr[0,6] = -4012
$ws.Range("H86:N86").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1791.4286
$arr[0,1] = 1801.6
$arr[0,2] = 1766
$arr[0,3] = 9008
$arr[0,4] = 8830
$arr[0,5] = -3392
$arr[0,6] = -20062
$ws.Range("H89:N89").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 925.86664
$arr[0,1] = 1107.091
$arr[0,2] = 427.5
$arr[0,3] = 1107.091
$arr[0,4] = 427.5
$arr[0,5] = -656.0909999999999
$arr[0,6] = -1329.5
$ws.Range("H94:N94").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10646.258
$arr[0,1] = 6260.0454
$arr[0,2] = 21368.111
$arr[0,3] = 6260.0454
$arr[0,4] = 21368.111
$arr[0,5] = -4340.0454
$arr[0,6] = -25208.111
$ws.Range("H107:N107").Value = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3281
$arr[0,1] = 4005
$arr[0,2] = 3100
$arr[0,3] = 4005
$arr[0,4] = 3100
$arr[0,5] = -3381
$arr[0,6] = -4348
$ws.Range("H62:N62").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3281
$arr[0,1] = 4005
$arr[0,2] = 3100
$arr[0,3] = 20025
$arr[0,4] = 15500
$arr[0,5] = -16905
$arr[0,6] = -21740
$ws.Range("H65:N65").Value = $arr

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2227.59
$arr[0,1] = 717.775
$arr[0,2] = 3234.1333
$arr[0,3] = 2153.325
$arr[0,4] = 9702.3999
$arr[0,5] = -1342.325
$arr[0,6] = -11324.3999
$ws.Range("H68:N68").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2227.59
$arr[0,1] = 717.775
$arr[0,2] = 3234.1333
$arr[0,3] = 6459.974999999999
$arr[0,4] = 29107.1997
$arr[0,5] = -2403.974999999999
$arr[0,6] = -37219.1997
$ws.Range("H71:N71").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 980.9375
$arr[0,1] = 454.2973
$arr[0,2] = 2752.3635
$arr[0,3] = 1362.8919
$arr[0,4] = 8257.0905
$arr[0,5] = 557.1080999999999
$arr[0,6] = -12097.0905
$ws.Range("H107:N107").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 29414576
$arr[0,1] = 3040
$arr[0,2] = 71431060
$arr[0,3] = 9120
$arr[0,4] = 214293180
$arr[0,5] = -4020
$arr[0,6] = -214303380
$ws.Range("H137:N137").Value = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 28439.83
$arr[0,1] = 37285.266
$arr[0,2] = 4315.909
$arr[0,3] = 37285.266
$arr[0,4] = 4315.909
$arr[0,5] = -37015.266
$arr[0,6] = -4855.909
$ws.Range("H70:N70").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 28439.83
$arr[0,1] = 37285.266
$arr[0,2] = 4315.909
$arr[0,3] = 37285.266
$arr[0,4] = 4315.909
$arr[0,5] = -36349.266
$arr[0,6] = -6187.909
$ws.Range("H73:N73").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 83334340
$arr[0,1] = 55556550
$arr[0,2] = 166667680
$arr[0,3] = 55556550
$arr[0,4] = 166667680
$arr[0,5] = -55556054
$arr[0,6] = -166668672
$ws.Range("H97:N97").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 74078780
$arr[0,1] = 153847090
$arr[0,2] = 8221.786
$arr[0,3] = 461541270
$arr[0,4] = 24665.358
$arr[0,5] = -461538740
$arr[0,6] = -29725.358
$ws.Range("H132:N132").Value = $arr

$ws.Range("H133").Value = 21900
$ws.Range("J133").Value = 21900
$ws.Range("L133").Value = 21900
$ws.Range("N133").Value = -32020

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 32000
$ws.Range("J137").Value = 32000
$ws.Range("L137").Value = 32000
$ws.Range("N137").Value = -42200

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 49333.332
$ws.Range("J140").Value = 49333.332
$ws.Range("L140").Value = 49333.332
$ws.Range("N140").Value = -59693.332

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1594.8718
$arr[0,1] = 1502.9412
$arr[0,2] = 2220
$arr[0,3] = 1502.9412
$arr[0,4] = 2220
$arr[0,5] = -753.9412
$arr[0,6] = -3718
$ws.Range("H68:N68").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1594.8718
$arr[0,1] = 1502.9412
$arr[0,2] = 2220
$arr[0,3] = 7514.706
$arr[0,4] = 11100
$arr[0,5] = -3770.706
$arr[0,6] = -18588
$ws.Range("H71:N71").Value = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10002290
$arr[0,1] = 1490.375
$arr[0,2] = 27781488
$arr[0,3] = 4471.125
$arr[0,4] = 83344464
$arr[0,5] = -2021.125
$arr[0,6] = -83349364
$ws.Range("H122:N122").Value = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6756.75
$arr[0,1] = 8864.5
$arr[0,2] = 433.5
$arr[0,3] = 26593.5
$arr[0,4] = 1300.5
$arr[0,5] = -24143.5
$arr[0,6] = -6200.5
$ws.Range("H122:N122").Value = $arr
